# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund-holdings detail, same shape as the
# other quarterly sheets) right before the "总计" (totals) sheet, and adds a
# matching summary row at the top of "总计".
#
# To land the new sheet's workbook-level sheetId at 6 (with 总计 shifting to
# 7, which is also how the underlying worksheet parts get renumbered -
# sheet6.xml/sheet7.xml) we delete "总计" first (freeing its sheetId), then
# re-create the two sheets in the right order, then repopulate "总计"'s data
# (shifted down by one row to make room for the new 2022-Q1 summary line).

$wb = $excel.ActiveWorkbook

# Style-donor cells already present in the workbook, used purely so the new
# cells pick up the same bold/centered "header" look (style index 2) that
# every other sheet in this workbook uses for its header row + leading
# index column.
$styleSrcSheet = $wb.Worksheets.Item("2021-Q4")
$headerStyleCell = $styleSrcSheet.Range("B1")
$indexStyleCell = $styleSrcSheet.Range("A2")

# --- 1. Re-create "总计" after a new "2022-Q1" sheet -----------------------

$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete()

$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4)
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

# --- 2. Populate "2022-Q1" (fund holdings detail) --------------------------

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$headerStyleCell.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1Rows = @(
    @(0, "501203", "易方达创新未来18个月封闭运作混合A", "76.88", "79.51", "2.96", "2.2756", 8),
    @(1, "110013", "易方达科翔混合", "49.44", "87.59", "3.77", "1.8639", 7),
    @(2, "009341", "易方达均衡成长股票", "64.76", "82.97", "2.54", "1.6449", 10),
    @(3, "110029", "易方达科讯混合", "36.09", "91.34", "3.09", "1.1152", 7),
    @(4, "009805", "国泰医药健康股票A", "12.14", "92.49", "8.87", "1.0768", 3),
    @(5, "110001", "易方达平稳增长混合", "33.39", "60.98", "2.98", "0.9950", 6),
    @(6, "160212", "国泰估值优势混合 (LOF)", "8.98", "62.69", "9.10", "0.8172", 1),
    @(7, "011326", "国泰医药健康股票C", "1.09", "92.49", "8.87", "0.0967", 3)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = "'" + $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = "'" + $row[3]
    $q1.Cells.Item($r, 5).Value = "'" + $row[4]
    $q1.Cells.Item($r, 6).Value = "'" + $row[5]
    $q1.Cells.Item($r, 7).Value = "'" + $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

$indexStyleCell.Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)

# --- 3. Populate "总计" (2022-Q1 row inserted, rest shifted down one) ------

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$headerStyleCell.Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$totalRows = @(
    @(0, "2022-Q1", 8, 9.890000000000001),
    @(1, "2021-Q4", 16, 15.88),
    @(2, "2021-Q3", 7, 8.76),
    @(3, "2021-Q2", 11, 9.1),
    @(4, "2021-Q1", 2, 0.46),
    @(5, "2020-Q4", 2, 0.41)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$indexStyleCell.Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
